# Add fun to extract titles from links
# Removes the two now-unused sheets, builds a "search word / link / link text"
# table with real hyperlinks on Sheet1, and formats the header row + link column.

$wb = $excel.ActiveWorkbook

# ---- drop the unused sheets -------------------------------------------------
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

$ws = $wb.Worksheets.Item(1)

# ---- column widths -----------------------------------------------------------
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 49.88

# ---- header row ---------------------------------------------------------------
$ws.Range("A1").Value = "Search Word"
$ws.Range("B1").Value = "search link"
$ws.Range("C1").Value = "Link Text"

$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# ---- search word column (A) ----------------------------------------------------
$searchWord = "260 متخصصًا عالميًا يترجمون أسئلة أولمبياد الكيمياء الدولي إلى 50 لغة"
$ws.Range("A2").Value = $searchWord
$ws.Range("A3").Value = $searchWord
$ws.Range("A4").Value = $searchWord
$ws.Range("A5").Value = $searchWord
$ws.Range("A6").Value = $searchWord
$ws.Range("A7").Value = $searchWord
$ws.Range("A8").Value = $searchWord
$ws.Range("A9").Value = $searchWord
$ws.Range("A10").Value = $searchWord

# ---- link-text column (C) ------------------------------------------------------
$ws.Range("C2").Value = '260 متخصصًا عالميًا يترجمون أسئلة أولمبياد الكيمياء الدولي إلى 50 لغةsabq.org › محليات'
$ws.Range("C3").Value = '260 متخصصاً عالمياً يترجمون أسئلة أولمبياد الكيمياء الدولي إلى 50 لغةlocal ‹ news ‹ www.okaz.com.sa'
$ws.Range("C4").Value = '260 متخصصًا عالميًا يترجمون أسئلة أولمبياد الكيمياء الدولي إلى 50 لغةأخبار وتقارير ‹ www.tech-mag.net'
$ws.Range("C5").Value = 'قلل الاستهلاك وحسن الأداء.. نصيحتان من "المواصفات" عن "الاستخدام ...post ‹ al-hadath.com'
$ws.Range("C6").Value = 'لحراطين - موقع تلماس الإخباري... ‹ www.tilmass.info'
$ws.Range("C7").Value = '"تعليم عسير" تحصل على 3 جوائز في الأولمبياد الوطني للإبداع بمجال ...... ‹ aldira.net'
$ws.Range("C8").Value = $searchWord
$ws.Range("C9").Value = $searchWord
$ws.Range("C10").Value = '260 متخصصًا عالميًا يترجمون أسئلة أولمبياد الكيمياء الدولي إلى 50 لغة ...'

# ---- link column (B) : real hyperlinks, display text = the URL itself ----------
$ws.Hyperlinks.Add($ws.Range("B2"), "https://sabq.org/saudia/vhlt14pdwv")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://www.okaz.com.sa/news/local/2166479")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.tech-mag.net/260-%D9%85%D8%AA%D8%AE%D8%B5%D8%B5%D9%8B%D8%A7-%D8%B9%D8%A7%D9%84%D9%85%D9%8A%D9%8B%D8%A7-%D9%8A%D8%AA%D8%B1%D8%AC%D9%85%D9%88%D9%86-%D8%A3%D8%B3%D8%A6%D9%84%D8%A9-%D8%A3%D9%88%D9%84%D9%85%D8%A8%D9%8A/")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://al-hadath.com/post/75279")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://www.tilmass.info/?tag=%D9%84%D8%AD%D8%B1%D8%A7%D8%B7%D9%8A%D9%86")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://aldira.net/?p=26334")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://sabq.org/saudia/vhlt14pdwv")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://shababeks.com/2024/07/23/mr-129/")
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.tech-mag.net/260-%D9%85%D8%AA%D8%AE%D8%B5%D8%B5%D9%8B%D8%A7-%D8%B9%D8%A7%D9%84%D9%85%D9%8A%D9%8B%D8%A7-%D9%8A%D8%AA%D8%B1%D8%AC%D9%85%D9%88%D9%86-%D8%A3%D8%B3%D8%A6%D9%84%D8%A9-%D8%A3%D9%88%D9%84%D9%85%D8%A8%D9%8A/")

# leftover hyperlink-styled (but empty) cells below the table
$ws.Range("B13").Style = "Hyperlink"
$ws.Range("B14").Style = "Hyperlink"
$ws.Range("B15").Style = "Hyperlink"

# ---- clean up left-over formatting from the old layout -------------------------
$ws.Range("A18").Clear()

# ---- selection / view state -----------------------------------------------------
$ws.Range("A4:XFD10").Select()
